# Commit: "added title parameter for bar plot"
#
# 1) Add a title ("Pearson's correlation") to the value axis of the
#    clustered-bar chart that lives on slide 12 ("Mean result for 1:3").
# 2) Add a small "channels" textbox under the chart on the same slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

# --- 1. Chart value-axis title -------------------------------------------
$chartShape = $s.Shapes.Item(2)
$chart = $chartShape.Chart
$valAx = $chart.Axes(2)          # 2 = xlValue
$valAx.HasTitle = $true
$axisTitle = $valAx.AxisTitle
$axisTitle.Text = "Pearson’s correlation"
# Keep the title rotated top-to-bottom along the value axis, matching the
# existing tick-label orientation on this chart.
$axisTitle.Orientation = 90
$axisTitle.Format.TextFrame2.TextRange.Font.Size = 10

# --- 2. "channels" textbox below the chart --------------------------------
$emuPerPt = 914400 / 72
$left   = 709864 / $emuPerPt
$top    = 5293895 / $emuPerPt
$width  = 818147 / $emuPerPt
$height = 230832 / $emuPerPt

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 2"
$tb.TextFrame.WordWrap = $true
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = $false
$tb.TextFrame.TextRange.Text = "channels"
$tb.TextFrame.TextRange.Font.Size = 9
